$wb = $excel.ActiveWorkbook
$ws2021 = $wb.Worksheets.Item(1)

# --- Insert a new sheet "2022-Q1" right after "2021-Q4" ---
$wsNew = $wb.Worksheets.Add($null, $ws2021)
$wsNew.Name = "2022-Q1"

# Header row: reuse the bordered/bold header style (s=2) from "2021-Q4",
# then overwrite the text (columns differ: scale / position / ranking are new).
$ws2021.Range("B1:H1").Copy()
$wsNew.Range("B1:H1").PasteSpecial(-4122)
$wsNew.Range("B1").Value = "基金代码"
$wsNew.Range("C1").Value = "基金名称"
$wsNew.Range("D1").Value = "基金规模"
$wsNew.Range("E1").Value = "股票总仓位"
$wsNew.Range("F1").Value = "仓位占比"
$wsNew.Range("G1").Value = "持有市值(亿元)"
$wsNew.Range("H1").Value = "仓位排名"

# Index column (A2:A10): same style as "2021-Q4"'s A-column (s=2)
$ws2021.Range("A2").Copy()
$wsNew.Range("A2:A10").PasteSpecial(-4122)

# Fund rows
$wsNew.Range("A2").Value = 0
$wsNew.Range("B2").Value = "'002340"
$wsNew.Range("C2").Value = "富国价值优势混合"
$wsNew.Range("D2").Value = "'69.56"
$wsNew.Range("E2").Value = "'93.66"
$wsNew.Range("F2").Value = "'2.63"
$wsNew.Range("G2").Value = "'1.8294"
$wsNew.Range("H2").Value = 10

$wsNew.Range("A3").Value = 1
$wsNew.Range("B3").Value = "'000729"
$wsNew.Range("C3").Value = "建信中小盘先锋股票"
$wsNew.Range("D3").Value = "'33.97"
$wsNew.Range("E3").Value = "'89.28"
$wsNew.Range("F3").Value = "'3.93"
$wsNew.Range("G3").Value = "'1.3350"
$wsNew.Range("H3").Value = 7

$wsNew.Range("A4").Value = 2
$wsNew.Range("B4").Value = "'005368"
$wsNew.Range("C4").Value = "富国清洁能源产业灵活配置混合A"
$wsNew.Range("D4").Value = "'23.52"
$wsNew.Range("E4").Value = "'88.60"
$wsNew.Range("F4").Value = "'4.71"
$wsNew.Range("G4").Value = "'1.1078"
$wsNew.Range("H4").Value = 5

$wsNew.Range("A5").Value = 3
$wsNew.Range("B5").Value = "'530005"
$wsNew.Range("C5").Value = "建信优化配置混合"
$wsNew.Range("D5").Value = "'21.57"
$wsNew.Range("E5").Value = "'86.97"
$wsNew.Range("F5").Value = "'3.95"
$wsNew.Range("G5").Value = "'0.8520"
$wsNew.Range("H5").Value = 7

$wsNew.Range("A6").Value = 4
$wsNew.Range("B6").Value = "'009693"
$wsNew.Range("C6").Value = "富国积极成长一年定期开放混合"
$wsNew.Range("D6").Value = "'17.82"
$wsNew.Range("E6").Value = "'98.74"
$wsNew.Range("F6").Value = "'3.33"
$wsNew.Range("G6").Value = "'0.5934"
$wsNew.Range("H6").Value = 7

$wsNew.Range("A7").Value = 5
$wsNew.Range("B7").Value = "'004674"
$wsNew.Range("C7").Value = "富国新机遇灵活配置混合A"
$wsNew.Range("D7").Value = "'23.28"
$wsNew.Range("E7").Value = "'93.61"
$wsNew.Range("F7").Value = "'2.41"
$wsNew.Range("G7").Value = "'0.5610"
$wsNew.Range("H7").Value = 9

$wsNew.Range("A8").Value = 6
$wsNew.Range("B8").Value = "'000756"
$wsNew.Range("C8").Value = "建信潜力新蓝筹股票"
$wsNew.Range("D8").Value = "'10.47"
$wsNew.Range("E8").Value = "'84.61"
$wsNew.Range("F8").Value = "'3.93"
$wsNew.Range("G8").Value = "'0.4115"
$wsNew.Range("H8").Value = 7

$wsNew.Range("A9").Value = 7
$wsNew.Range("B9").Value = "'004675"
$wsNew.Range("C9").Value = "富国新机遇灵活配置混合C"
$wsNew.Range("D9").Value = "'3.84"
$wsNew.Range("E9").Value = "'93.61"
$wsNew.Range("F9").Value = "'2.41"
$wsNew.Range("G9").Value = "'0.0925"
$wsNew.Range("H9").Value = 9

$wsNew.Range("A10").Value = 8
$wsNew.Range("B10").Value = "'011127"
$wsNew.Range("C10").Value = "富国清洁能源产业灵活配置混合C"
$wsNew.Range("D10").Value = "'1.61"
$wsNew.Range("E10").Value = "'88.60"
$wsNew.Range("F10").Value = "'4.71"
$wsNew.Range("G10").Value = "'0.0758"
$wsNew.Range("H10").Value = 5

# Columns B,D,E,F,G hold numeric-looking text (fund codes / percentages);
# clear the auto "quote prefix" styling picked up from the apostrophe so the
# cells stay plain (unstyled) text, matching the source data.
$wsNew.Range("B2:B10").Style = "Normal"
$wsNew.Range("D2:G10").Style = "Normal"

# --- "总计" sheet: insert a new top data row for 2022-Q1, pushing 2021-Q4 down ---
$wsTotal = $wb.Worksheets.Item(3)
$wsTotal.Rows.Item(2).Insert()

$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)
$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 9
$wsTotal.Range("D2").Value = 6.86
$wsTotal.Range("B2:D2").Style = "Normal"

Write-Output "done"
